$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4960103631019592
$ws.Range("B1").Value = 1.315841436386108
$ws.Range("C1").Value = 6.349945068359375
$ws.Range("D1").Value = 1.672683238983154
$ws.Range("E1").Value = 1.522847414016724
